$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value (Price and Volume(1h) columns)
$changes = @{
    'D2' = '71.272.39'
    'E2' = '  +1.41%  '
    'D3' = '3.865.68'
    'E3' = '  +1.72%  '
    'E4' = '  +0.12%  '
    'D5' = '692.17'
    'E5' = '  +3.29%  '
    'D6' = '173.22'
    'E6' = '  +2.74%  '
    'D7' = '3.863.54'
    'E7' = '  +1.66%  '
    'E9' = '  +0.20%  '
    'E10' = '  +1.88%  '
    'D11' = '7.36'
    'E11' = '  +4.09%  '
    'D12' = '0.463'
    'E12' = '  +0.34%  '
    'E13' = '  +6.50%  '
    'D14' = '36.64'
    'E14' = '  +2.79%  '
    'D15' = '4.504.63'
    'E15' = '  +1.51%  '
    'D16' = '3.906.03'
    'E16' = '  +3.34%  '
    'D17' = '71.314.49'
    'E17' = '  +1.50%  '
    'D18' = '17.81'
    'E18' = '  +0.87%  '
    'E19' = '  +1.20%  '
    'E20' = '  +0.29%  '
    'D21' = '11.09'
    'E21' = '  -2.92%  '
    'D22' = '494.81'
    'E22' = '  +4.18%  '
    'D23' = '0.724'
    'E23' = '  +1.61%  '
    'D24' = '84.95'
    'E24' = '  +1.97%  '
    'E25' = '  +4.14%  '
    'D26' = '12.38'
    'E26' = '  +1.64%  '
    'D27' = '10.62'
    'E27' = '  +3.47%  '
    'D28' = '2.15'
    'E28' = '  +2.11%  '
    'D29' = '4.024.32'
    'E29' = '  +1.92%  '
    'E30' = '  +0.04%  '
    'E31' = '  +10.32%  '
    'E32' = '  +3.67%  '
    'E33' = '  +0.24%  '
    'D34' = '29.82'
    'E34' = '  +0.88%  '
    'D35' = '0.179'
    'E35' = '  +0.34%  '
    'D36' = '9.33'
    'E36' = '  +2.56%  '
    'D37' = '3.817.19'
    'E37' = '  +1.66%  '
    'D38' = '0.999'
    'E38' = '  -0.08%  '
    'E39' = '  +2.68%  '
    'D40' = '2.39'
    'E40' = '  +13.02%  '
    'D41' = '3.44'
    'E41' = '  +1.92%  '
    'E42' = '  +2.03%  '
    'E43' = '  +6.55%  '
    'E44' = '  +0.06%  '
    'E45' = '  +0.06%  '
    'D46' = '164.20'
    'E46' = '  +3.18%  '
    'E47' = '  +5.66%  '
    'D48' = '48.67'
    'E48' = '  +1.33%  '
    'D49' = '44.50'
    'E49' = '  -2.01%  '
    'E50' = '  +1.38%  '
    'D51' = '8.70'
    'E51' = '  +2.38%  '
}

foreach ($addr in $changes.Keys) {
    $cell = $ws.Range($addr)
    # Force text storage so numeric-looking strings (e.g. "692.17") keep
    # their exact original text formatting instead of being coerced to a
    # floating point number by Excel's automatic type detection.
    $cell.NumberFormat = "@"
    $cell.Value = $changes[$addr]
}
